# Updates cryptos list values (Price / Volume(1h) columns, plus the
# BabyDogeCoin insertion that bumps TheSandbox down a row and drops Aptos)
# per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.402.35'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.876.43'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '0.7123'
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").Value = '241.82'
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.3113'
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").Value = '0.07716'
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("D10").Value = '25.39'
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = '0.08378'
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("D12").Value = '1.886.99'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = "'5.250"
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '0.7166'
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").Value = '29.416.34'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").Value = '0.000008247'
$ws.Range("E17").Value = '  +5.35%  '
$ws.Range("D18").Value = '5.978'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("D19").Value = '244.09'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '2.133.53'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '0.9994'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '7.899'
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '0.1617'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '163.78'
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("D27").Value = '9.027'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").Value = '18.59'
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = '1.507'
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("D30").Value = '4.415'
$ws.Range("E30").Value = '  +0.83%  '
$ws.Range("D31").Value = '1.293'
$ws.Range("E31").Value = '  -4.14%  '
$ws.Range("D32").Value = '4.332'
$ws.Range("E32").Value = '  +5.55%  '
$ws.Range("D33").Value = '0.05224'
$ws.Range("E33").Value = '  +0.69%  '
$ws.Range("D34").Value = '1.927'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("D35").Value = '0.7761'
$ws.Range("E35").Value = '  +6.88%  '
$ws.Range("D36").Value = '1.176'
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").Value = '2.681'
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").Value = '0.01867'
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("D40").Value = '1.165.05'
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("D41").Value = '6.411'
$ws.Range("E41").Value = '  +4.56%  '
$ws.Range("D42").Value = '73.48'
$ws.Range("E42").Value = '  +0.97%  '
$ws.Range("D43").Value = '0.8911'
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("D44").Value = '104.58'
$ws.Range("E44").Value = '  +2.53%  '
$ws.Range("D45").Value = '0.9999'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '2.032.68'
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").Value = '1.797'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").Value = '9.394'
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = "'0.00000000120"
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '0.4313'
$ws.Range("E51").Value = '  +0.66%  '
